$wb = $excel.ActiveWorkbook

# Sheet ALC row 4 (item 5470)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 167.5
$ws.Range("I4").Value = 167.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 167.5
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

# Sheet ALC row 9 (item 5487)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 134.75
$ws.Range("I9").Value = 155.6
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 155.6
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = 13.40000000000001
$ws.Range("N9").Value = -438

# Sheet ALC row 112 (item 27960)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1192.57
$ws.Range("I112").Value = 990
$ws.Range("J112").Value = 1196.7041
$ws.Range("K112").Value = 2970
$ws.Range("L112").Value = 3590.1123
$ws.Range("M112").Value = -1862
$ws.Range("N112").Value = -5806.1123

# Sheet ALC row 137 (item 44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2978409.5
$ws.Range("I137").Value = 6411999
$ws.Range("J137").Value = 2632.1333
$ws.Range("K137").Value = 19235997
$ws.Range("L137").Value = 7896.3999
$ws.Range("M137").Value = -19233447
$ws.Range("N137").Value = -12996.3999

# Sheet ALC row 138 (item 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4503.746
$ws.Range("I138").Value = 5088.154
$ws.Range("J138").Value = 4363.0557
$ws.Range("K138").Value = 15264.462
$ws.Range("L138").Value = 13089.1671
$ws.Range("M138").Value = -10124.462
$ws.Range("N138").Value = -23369.1671

# Sheet ARM row 61 (item 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 13891540
$ws.Range("I61").Value = 23811640
$ws.Range("J61").Value = 3400
$ws.Range("K61").Value = 23811640
$ws.Range("L61").Value = 3400
$ws.Range("M61").Value = -23811428
$ws.Range("N61").Value = -3824

# Sheet ARM row 63 (item 12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 33743.098
$ws.Range("I63").Value = 302666.66
$ws.Range("J63").Value = 4929.857
$ws.Range("K63").Value = 302666.66
$ws.Range("L63").Value = 4929.857
$ws.Range("M63").Value = -301980.66
$ws.Range("N63").Value = -6301.857

# Sheet ARM row 66 (item 12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 33743.098
$ws.Range("I66").Value = 302666.66
$ws.Range("J66").Value = 4929.857
$ws.Range("K66").Value = 1513333.3
$ws.Range("L66").Value = 24649.285
$ws.Range("M66").Value = -1509901.3
$ws.Range("N66").Value = -31513.285

# Sheet ARM row 74 (item 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9092337
$ws.Range("I74").Value = 830.8205
$ws.Range("J74").Value = 31252884
$ws.Range("K74").Value = 830.8205
$ws.Range("L74").Value = 31252884
$ws.Range("M74").Value = 43.17949999999996
$ws.Range("N74").Value = -31254632

# Sheet ARM row 77 (item 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 9092337
$ws.Range("I77").Value = 830.8205
$ws.Range("J77").Value = 31252884
$ws.Range("K77").Value = 4154.1025
$ws.Range("L77").Value = 156264420
$ws.Range("M77").Value = 213.8975
$ws.Range("N77").Value = -156273156

# Sheet ARM row 132 (item 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2335642
$ws.Range("I132").Value = 3680
$ws.Range("J132").Value = 5923275.5
$ws.Range("K132").Value = 11040
$ws.Range("L132").Value = 17769826.5
$ws.Range("M132").Value = -8510
$ws.Range("N132").Value = -17774886.5

# Sheet ARM row 136 (item 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 13891540
$ws.Range("I136").Value = 23811640
$ws.Range("J136").Value = 3400
$ws.Range("K136").Value = 71434920
$ws.Range("L136").Value = 10200
$ws.Range("M136").Value = -71432370
$ws.Range("N136").Value = -15300

# Sheet BSM row 82 (item 11877)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 12886.267
$ws.Range("I82").Value = 5757.2856
$ws.Range("J82").Value = 19124.125
$ws.Range("K82").Value = 5757.2856
$ws.Range("L82").Value = 19124.125
$ws.Range("M82").Value = -5374.2856
$ws.Range("N82").Value = -19890.125

# Sheet BSM row 85 (item 11877)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 12886.267
$ws.Range("I85").Value = 5757.2856
$ws.Range("J85").Value = 19124.125
$ws.Range("K85").Value = 5757.2856
$ws.Range("L85").Value = 19124.125
$ws.Range("M85").Value = -4431.2856
$ws.Range("N85").Value = -21776.125

# Sheet BSM row 86 (item 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1980.875
$ws.Range("I86").Value = 1979.1364
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 1979.1364
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -856.1364000000001
$ws.Range("N86").Value = -4246

# Sheet BSM row 89 (item 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1980.875
$ws.Range("I89").Value = 1979.1364
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 9895.682000000001
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -4279.682000000001
$ws.Range("N89").Value = -21232

# Sheet BSM row 134 (item 43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2947.2092
$ws.Range("I134").Value = 2980.2942
$ws.Range("J134").Value = 2822.2222
$ws.Range("K134").Value = 8940.882599999999
$ws.Range("L134").Value = 8466.6666
$ws.Range("M134").Value = -6405.882599999999
$ws.Range("N134").Value = -13536.6666

# Sheet CRP row 7 (item 5361)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 72.2
$ws.Range("I7").Value = 40.166668
$ws.Range("J7").Value = 120.25
$ws.Range("K7").Value = 40.166668
$ws.Range("L7").Value = 120.25
$ws.Range("M7").Value = 72.833332
$ws.Range("N7").Value = -346.25

# Sheet CRP row 31 (item 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7968.0557
$ws.Range("I31").Value = 4076.875
$ws.Range("J31").Value = 9079.821
$ws.Range("K31").Value = 4076.875
$ws.Range("L31").Value = 9079.821
$ws.Range("M31").Value = -3781.875
$ws.Range("N31").Value = -9669.821

# Sheet CRP row 34 (item 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7968.0557
$ws.Range("I34").Value = 4076.875
$ws.Range("J34").Value = 9079.821
$ws.Range("K34").Value = 4076.875
$ws.Range("L34").Value = 9079.821
$ws.Range("M34").Value = -3874.875
$ws.Range("N34").Value = -9483.821

# Sheet CRP row 105 (item 19928)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 768.5714
$ws.Range("I105").Value = 645
$ws.Range("J105").Value = 933.3333
$ws.Range("K105").Value = 645
$ws.Range("L105").Value = 933.3333
$ws.Range("M105").Value = 1102
$ws.Range("N105").Value = -4427.3333

# Sheet CRP row 122 (item 36196)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2057.3713
$ws.Range("I122").Value = 1767.0476
$ws.Range("J122").Value = 2492.8572
$ws.Range("K122").Value = 5301.142800000001
$ws.Range("L122").Value = 7478.571599999999
$ws.Range("M122").Value = -2851.142800000001
$ws.Range("N122").Value = -12378.5716

# Sheet CRP row 132 (item 44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 11301084
$ws.Range("I132").Value = 14707291
$ws.Range("J132").Value = 6668643
$ws.Range("K132").Value = 44121873
$ws.Range("L132").Value = 20005929
$ws.Range("M132").Value = -44119343
$ws.Range("N132").Value = -20010989

# Sheet CRP row 134 (item 44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4905670.5
$ws.Range("I134").Value = 6253821
$ws.Range("J134").Value = 3304.182
$ws.Range("K134").Value = 18761463
$ws.Range("L134").Value = 9912.545999999998
$ws.Range("M134").Value = -18758928
$ws.Range("N134").Value = -14982.546

# Sheet CUL row 34 (item 4749)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 7463246
$ws.Range("I34").Value = 246.83333
$ws.Range("J34").Value = 8197311.5
$ws.Range("K34").Value = 740.49999
$ws.Range("L34").Value = 24591934.5
$ws.Range("M34").Value = -656.49999
$ws.Range("N34").Value = -24592102.5

# Sheet CUL row 39 (item 4712)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2154.8823
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 2154.8823
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 6464.646900000001
$ws.Range("N39").Value = -7052.646900000001

# Sheet CUL row 107 (item 27838)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1423.7333
$ws.Range("I107").Value = 209.09091
$ws.Range("J107").Value = 2126.9473
$ws.Range("K107").Value = 627.27273
$ws.Range("L107").Value = 6380.841899999999
$ws.Range("M107").Value = 1292.72727
$ws.Range("N107").Value = -10220.8419

# Sheet GSM row 2 (item 5062)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 72.71429000000001
$ws.Range("I2").Value = 72.71429000000001
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 72.71429000000001
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

# Sheet GSM row 57 (item 2876)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 19999
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 19999
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 19999
$ws.Range("N57").Value = -21639

# Sheet GSM row 132 (item 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 35720690
$ws.Range("I132").Value = 62508524
$ws.Range("J132").Value = 3568.6667
$ws.Range("K132").Value = 187525572
$ws.Range("L132").Value = 10706.0001
$ws.Range("M132").Value = -187523042
$ws.Range("N132").Value = -15766.0001

# Sheet GSM row 140 (item 42458)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 54250
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 54250
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 54250
$ws.Range("N140").Value = -64610

# Sheet LTW row 16 (item 5289)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1100.9412
$ws.Range("I16").Value = 324.46155
$ws.Range("J16").Value = 3624.5
$ws.Range("K16").Value = 324.46155
$ws.Range("L16").Value = 3624.5
$ws.Range("M16").Value = -154.46155
$ws.Range("N16").Value = -3964.5

# Sheet LTW row 136 (item 44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3877165.5
$ws.Range("I136").Value = 1036.8975
$ws.Range("J136").Value = 41669420
$ws.Range("K136").Value = 3110.6925
$ws.Range("L136").Value = 125008260
$ws.Range("M136").Value = -560.6925000000001
$ws.Range("N136").Value = -125013360

# Sheet WVR row 93 (item 19613)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 33697.25
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 33697.25
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 33697.25
$ws.Range("N93").Value = -38689.25

# Sheet WVR row 132 (item 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3789802
$ws.Range("I132").Value = 1849.4
$ws.Range("J132").Value = 17159046
$ws.Range("K132").Value = 5548.200000000001
$ws.Range("L132").Value = 51477138
$ws.Range("M132").Value = -3018.200000000001
$ws.Range("N132").Value = -51482198

# Sheet WVR row 136 (item 44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1951.638
$ws.Range("I136").Value = 1796.4565
$ws.Range("J136").Value = 2546.5
$ws.Range("K136").Value = 5389.3695
$ws.Range("L136").Value = 7639.5
$ws.Range("M136").Value = -2839.3695
$ws.Range("N136").Value = -12739.5
